# "Module Code Generator" workbook update:
# CHAMP Program Updated in AOD, MH, MHBC, MHLC, and Hotlines
#
# The sheet works like a rolling log: row 2 holds the inputs for the
# resource currently being coded, rows 4-6 compute the generated R
# snippets from row 2 via formulas, and rows 9-15 keep a short history
# of previously generated snippets (plus some retired lookup data in
# row 12). Each "run" pushes the previous rows 9-11 contents down into
# 13-15 and drops a fresh result into 9-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the working input row (row 2) with the new resource -
#    MHLC_CHAMP / CHAMP Program (OASAS) - which drives the row 4-6
#    formulas to recompute automatically. B2/D2 pick up Text format to
#    match C2 (mirrors the source workbook's cell styling change).
$ws.Range("A2").Value = "MHLC_CHAMP"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "CHAMP Program (OASAS)"
$ws.Range("C2").Value = "MHLC()"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "CHAMP Program (OASAS)"

# 2. Push the freshly (re)computed row 4-6 snapshot (MHLC_CHAMP, now
#    reflecting the updated row-2 inputs) down into the history rows
#    13-15.
$ws.Range("A13").Value = $ws.Range("B4").Value2
$ws.Range("A14").Value = $ws.Range("B5").Value2
$ws.Range("A15").Value = $ws.Range("B6").Value2

# 3. The previous history block (rows 9-12) is retired: clear the old
#    lookup cells A12:D12 (A12/B12 fully cleared, C12/D12 just emptied
#    of their hyperlinked text but keep their formatting) and drop the
#    two hyperlinks that lived there.
$ws.Hyperlinks.Delete()
$ws.Range("A12").Clear()
$ws.Range("B12").Clear()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()

# 4. Drop in the newly generated MHBC_CHAMP snippets as the latest
#    history entry (rows 9-11).
$ws.Range("A9").Value = "mod_Accordion_ui('MHBC_CHAMP')"
$ws.Range("A10").Value = "mod_Accordion_server('MHBC_CHAMP', selector=selection, data=MHBC(), title = c('CHAMP Program (OASAS)'), Visible = T)"
$ws.Range("A11").Value = "mod_info_server('MHBC_CHAMP', selector = selection, data = MHBC(), rownametitle = c('CHAMP Program (OASAS)'), phone = T, website = T)"

# 5. Update the active selection / scroll position to match where the
#    user ended up after making the edits.
$ws.Range("B12").Select()
